$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51-71 down to 52-72.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with a new Papaya price entry.
# Descriptive columns mirror the surrounding "Vega Modelo de Temuco" / Papaya rows.
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = 'Vega Modelo de Temuco'
$ws.Range("C51").Value = 'La Araucanía'
$ws.Range("D51").Value = 44627
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 'Fruta'
$ws.Range("G51").Value = 100108
$ws.Range("H51").Value = 'Tropicales y subtropicales'
$ws.Range("I51").Value = 100108004
$ws.Range("J51").Value = 'Papaya'
$ws.Range("K51").Value = 'Cultivar IV Región'
$ws.Range("L51").Value = 'Primera'
$ws.Range("M51").Value = 80
$ws.Range("N51").Value = 22000
$ws.Range("O51").Value = 22000
$ws.Range("P51").Value = 22000
$ws.Range("Q51").Value = '$/bandeja 10 kilos'
$ws.Range("R51").Value = 'Provincia del Elquí'
$ws.Range("S51").Value = 2200
$ws.Range("T51").Value = 10
